# Mark the first three test cases (rows 2-4) as passed:
#   - "Actual Outcome" (column F) -> "Same as expected outcome."
#   - "Fail/Pass"       (column G) -> "Pass"
#
# Set the G (Pass) cells before the F (Same as expected outcome.) cells so
# that the shared-string table allocates "Pass" before "Same as expected
# outcome.", matching how Excel ordered the new entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Pass"
$ws.Range("G3").Value = "Pass"
$ws.Range("G4").Value = "Pass"

$ws.Range("F2").Value = "Same as expected outcome."
$ws.Range("F3").Value = "Same as expected outcome."
$ws.Range("F4").Value = "Same as expected outcome."

# Reflect the updated view: scrolled right so column C is left-most, with
# F2:F4 selected (active cell F2).
try {
    $excel.ActiveWindow.ScrollColumn = 3
} catch {
}

$ws.Range("F2:F4").Select() | Out-Null
